# This script applies a cyclic rotation of the full content of rows
# 55, 56, 58 and 59 (55 <- 59 <- 56 <- 58 <- 55) and a full swap of the
# content of rows 60 and 61, on the single worksheet "Artfynd".
#
# Whole rows (columns A:AY, matching the sheet's used range) are staged
# via Range.Copy into unused rows far below the data (1000+) so that the
# 4-way rotation does not clobber source data before it has been copied
# out. Range.Copy (rather than assigning .Value/.Value2 arrays) is used
# deliberately: it duplicates the literal stored cell content (numbers,
# booleans, text) without Excel's "smart" re-interpretation of
# date-looking strings (e.g. "2023-09-03") into date serial numbers,
# and it does not touch cell styles.
#
# After copying a staged row onto its destination row, any column that
# was *absent* (not just blank) in that row's original source is
# cleared with ClearContents() on the destination so the destination
# row ends up with exactly the same set of present/absent cells as the
# row it was copied from (Copy alone creates empty placeholder cells
# for every column in the copied range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"

function Get-RowRange($row) {
    return $ws.Range("A" + $row + ":" + $lastCol + $row)
}

# Columns (by letter) that are *not present at all* in the original XML
# for each source row (i.e. not even an empty cell). These are cleared
# on the destination after the row content has been copied there, so
# the destination matches the source's cell layout exactly.
$absentCols = @{
    55 = @("L","M","O","X","AC","AI","AL","AN","AP","AQ","AR","AS","AU","AV")
    56 = @("M","O","X","AC","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")
    58 = @("L","M","O","X","AC","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")
    59 = @("M","O","X","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")
    60 = @("M","O","X","AC","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")
    61 = @("M","O","X","AC","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AU","AV")
}

function Apply-Absent($destRow, $sourceRow) {
    foreach ($col in $absentCols[$sourceRow]) {
        $ws.Range($col + $destRow).ClearContents()
    }
}

# --- Stage the original content of the six affected rows into unused
#     rows (1000-1005), far below any real data, so the subsequent
#     writes don't overwrite data that is still needed. ---
$stage55 = 1000
$stage56 = 1001
$stage58 = 1002
$stage59 = 1003
$stage60 = 1004
$stage61 = 1005

(Get-RowRange 55).Copy((Get-RowRange $stage55))
(Get-RowRange 56).Copy((Get-RowRange $stage56))
(Get-RowRange 58).Copy((Get-RowRange $stage58))
(Get-RowRange 59).Copy((Get-RowRange $stage59))
(Get-RowRange 60).Copy((Get-RowRange $stage60))
(Get-RowRange 61).Copy((Get-RowRange $stage61))

# --- Write the rotated / swapped content back to the real rows. ---
# Row 55 <- original row 59
(Get-RowRange $stage59).Copy((Get-RowRange 55))
Apply-Absent 55 59

# Row 58 <- original row 55
(Get-RowRange $stage55).Copy((Get-RowRange 58))
Apply-Absent 58 55

# Row 56 <- original row 58
(Get-RowRange $stage58).Copy((Get-RowRange 56))
Apply-Absent 56 58

# Row 59 <- original row 56
(Get-RowRange $stage56).Copy((Get-RowRange 59))
Apply-Absent 59 56

# Row 60 <- original row 61
(Get-RowRange $stage61).Copy((Get-RowRange 60))
Apply-Absent 60 61

# Row 61 <- original row 60
(Get-RowRange $stage60).Copy((Get-RowRange 61))
Apply-Absent 61 60

# --- Clean up the staging rows so they don't leave any trace (and the
#     sheet's used range / dimension stays A1:AY61). ---
(Get-RowRange $stage55).ClearContents()
(Get-RowRange $stage56).ClearContents()
(Get-RowRange $stage58).ClearContents()
(Get-RowRange $stage59).ClearContents()
(Get-RowRange $stage60).ClearContents()
(Get-RowRange $stage61).ClearContents()

Write-Output "Row rotation/swap applied."
